# New daily price record for "Rabanito" in Vega Modelo de Temuco: a row is
# inserted at row 23 (pushing the existing rows 23-103 down to 24-104, which
# grows the used range from A1:R103 to A1:R104), and the new row 23 is filled
# in with that day's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(23).Insert()

$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 45054
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 300000001
$ws.Range("G23").Value = "Rabanito"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 8000
$ws.Range("N23").Value = "$/docena de paquetes"
$ws.Range("O23").Value = "Provincia de Cautín"
$ws.Range("P23").Value = 667
$ws.Range("Q23").Value = 12
$ws.Range("R23").Value = "Hortaliza"
